{"js": "// Update the title placeholder text, and remove the now-unused\n// UNCLASSIFIED / country / BLUF / body paragraphs that used to follow it,\n// leaving just the (retitled) title paragraph in the document body.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// 1) Rewrite the title placeholder text in-place (keeps the Title style).\nconst title = body.paragraphs.items[0];\ntitle.insertText(\"[Intelligence Note or Reporting Highlights]\", \"Replace\");\nawait context.sync();\n\n// 2) Drop every paragraph that follows the title. Office.js paragraph\n//    proxies are live references, so re-fetch the collection and always\n//    delete the paragraph right after the title (index 1) -- deleting from\n//    the back would get \"stuck\" on the document's final paragraph mark,\n//    which Word never lets you remove outright.\nwhile (true) {\n  body.paragraphs.load(\"items\");\n  await context.sync();\n  if (body.paragraphs.items.length <= 1) break;\n  body.paragraphs.items[1].delete();\n  await context.sync();\n}\n", "ps1": "# Update the title placeholder text, and remove the now-unused\n# UNCLASSIFIED / country / BLUF / body paragraphs that used to follow it,\n# leaving just the (retitled) title paragraph in the document body.\n\n$d = $word.ActiveDocument\n\n# 1) Rewrite the title placeholder text in-place (keeps the Title style,\n#    since only the run text inside the paragraph range is replaced).\n$d.Paragraphs.Item(1).Range.Text = \"[Intelligence Note or Reporting Highlights]\"\n\n# 2) Drop every paragraph that follows the title. Always delete the\n#    paragraph right after the title (index 2) -- deleting from the back\n#    would get \"stuck\" on the document's final paragraph mark, which Word\n#    never lets you remove outright.\nwhile ($d.Paragraphs.Count -gt 1) {\n    $d.Paragraphs.Item(2).Range.Delete()\n}\n"}
